$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) from 45174 (2023-09-05)
# to 45175 (2023-09-06) for rows 2 through 7.
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45175
}
